$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.791.22"
$ws.Range("E2").Value = "  +1.07%  "

$ws.Range("D3").Value = "3.848.84"
$ws.Range("E3").Value = "  +0.02%  "

$ws.Range("D5").Value = "'601.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.68%  "

$ws.Range("D6").Value = "'171.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.78%  "

$ws.Range("D7").Value = "3.847.43"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("E10").Value = "  +2.75%  "

$ws.Range("E11").Value = "  +3.32%  "

$ws.Range("E12").Value = "  +0.05%  "

$ws.Range("D13").Value = "'0.0000285"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +14.60%  "

$ws.Range("D14").Value = "'37.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.11%  "

$ws.Range("D15").Value = "4.496.72"
$ws.Range("E15").Value = "  +0.09%  "

$ws.Range("D16").Value = "3.928.02"
$ws.Range("E16").Value = "  +1.91%  "

$ws.Range("D17").Value = "68.784.46"
$ws.Range("E17").Value = "  +0.92%  "

$ws.Range("D18").Value = "'18.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("D19").Value = "'7.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.93%  "

$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").Value = "'11.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.42%  "

$ws.Range("D22").Value = "'475.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.29%  "

$ws.Range("D23").Value = "'0.728"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.49%  "

$ws.Range("D24").Value = "'0.0000164"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.43%  "

$ws.Range("D25").Value = "'83.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.53%  "

$ws.Range("E26").Value = "  +0.92%  "

$ws.Range("D27").Value = "'12.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.34%  "

$ws.Range("D28").Value = "'10.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.95%  "

$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").Value = "'2.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.64%  "

$ws.Range("D31").Value = "4.001.94"
$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("E33").Value = "  +1.14%  "

$ws.Range("E34").Value = "  -0.29%  "

$ws.Range("D35").Value = "'9.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.86%  "

$ws.Range("D36").Value = "3.816.61"
$ws.Range("E36").Value = "  -0.13%  "

$ws.Range("D37").Value = "'4.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +21.23%  "

$ws.Range("E38").Value = "  -0.51%  "

$ws.Range("D39").Value = "'0.141"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.62%  "

$ws.Range("E40").Value = "  +0.84%  "

$ws.Range("E41").Value = "  +0.33%  "

$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("E43").Value = "  +1.30%  "

$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("D45").Value = "'0.000302"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.78%  "

$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("D47").Value = "'419.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.90%  "

$ws.Range("D48").Value = "'8.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.08%  "

$ws.Range("D49").Value = "'46.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.04%  "

$ws.Range("D50").Value = "'142.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("E51").Value = "  -0.06%  "

